$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: E7 changes from "self01.jpg" to "adult.jpg"
$ws.Range("E7").Value = "adult.jpg"

# Copy the formatting (style) of row 7's populated cells onto the new rows 8 and 9
foreach ($col in @("A","B","C","D","E","H")) {
    $ws.Range($col + "7").Copy()
    $ws.Range($col + "8").PasteSpecial(-4122)
    $ws.Range($col + "9").PasteSpecial(-4122)
}

# Row 8 (id 7): "月光光"
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "月光光"
$ws.Range("C8").Value = "月光光"
$ws.Range("D8").Value = "00007.png"
$ws.Range("E8").Value = "fright.jpg"
$ws.Range("H8").Value = "中秋的大月光在图中的天空出現。月光佔了天空一半，又圆又大"

# Row 9 (id 8): "暴龍頭上"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "暴龍頭上"
$ws.Range("C9").Value = "暴龍頭上"
$ws.Range("D9").Value = "00008.png"
$ws.Range("E9").Value = "young.jpg"
$ws.Range("H9").Value = "图中的人坐在暴龍頭上，人的服裝换上樹葉做的和背景是朱羅记時代"
